$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 179
$ws.Range("F3").Value = 2452
$ws.Range("F6").Value = 1040
$ws.Range("F7").Value = 92
$ws.Range("F8").Value = 305
$ws.Range("F9").Value = 587
$ws.Range("F10").Value = 3450
$ws.Range("F11").Value = 911
$ws.Range("F12").Value = 1131
$ws.Range("F15").Value = 872
$ws.Range("F16").Value = 7
$ws.Range("F17").Value = 1127
$ws.Range("F18").Value = 1767
$ws.Range("F21").Value = 1529
$ws.Range("F22").Value = 1070
$ws.Range("F23").Value = 1017
$ws.Range("F25").Value = 4170
$ws.Range("F26").Value = 31
$ws.Range("F27").Value = 2674
$ws.Range("F28").Value = 1198

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 30
$ws.Range("F11").Value = 6
$ws.Range("F19").Value = 23
$ws.Range("F20").Value = 23
$ws.Range("F26").Value = 80
$ws.Range("F36").Value = 60
$ws.Range("F39").Value = 409
$ws.Range("F40").Value = 5

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 2525
$ws.Range("F6").Value = 2534
$ws.Range("F7").Value = 9566
$ws.Range("F8").Value = 152
$ws.Range("F12").Value = 2915
$ws.Range("F13").Value = 429
$ws.Range("F14").Value = 770

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2525
$ws.Range("F4").Value = 152
$ws.Range("F5").Value = 179
$ws.Range("F6").Value = 2452
$ws.Range("F8").Value = 2915
$ws.Range("F10").Value = 770
$ws.Range("F13").Value = 92
$ws.Range("F14").Value = 305
$ws.Range("F15").Value = 587
$ws.Range("F16").Value = 30
$ws.Range("F17").Value = 911
$ws.Range("F18").Value = 1131
$ws.Range("F21").Value = 872
$ws.Range("F24").Value = 1127
$ws.Range("F27").Value = 23
$ws.Range("F28").Value = 23
$ws.Range("F32").Value = 1767
$ws.Range("F35").Value = 1529
$ws.Range("F37").Value = 80
$ws.Range("F38").Value = 80
$ws.Range("F39").Value = 1070
$ws.Range("F41").Value = 1017
$ws.Range("F44").Value = 60
$ws.Range("F45").Value = 4170
$ws.Range("F46").Value = 409
$ws.Range("F47").Value = 2674
$ws.Range("F48").Value = 5
